$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.07846107257518
$ws.Range("C2").Value = 8.448949366929495
$ws.Range("D2").Value = 9.482794649332133
$ws.Range("E2").Value = 13.77982124675057
$ws.Range("F2").Value = 31.18773958003534
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 9.838890416809909
$ws.Range("N2").Value = 16.86215702231468
$ws.Range("O2").Value = 23.16503057641329
$ws.Range("B3").Value = 14.54113448007725
$ws.Range("C3").Value = 7.961981145490893
$ws.Range("D3").Value = 9.440132079982726
$ws.Range("E3").Value = 13.72411774595267
$ws.Range("F3").Value = 31.16627103514676
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 9.844916739155247
$ws.Range("N3").Value = 16.91142088730066
$ws.Range("O3").Value = 23.20166491549962
$ws.Range("B4").Value = 14.20295255604326
$ws.Range("C4").Value = 7.647229857235717
$ws.Range("D4").Value = 9.415314629859827
$ws.Range("E4").Value = 13.69266667445603
$ws.Range("F4").Value = 31.16229911810432
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 9.85027239352274
$ws.Range("N4").Value = 16.9435513603388
$ws.Range("O4").Value = 23.23064480849427
$ws.Range("B5").Value = 14.06327637199842
$ws.Range("C5").Value = 7.515091603985273
$ws.Range("D5").Value = 9.405554902380937
$ws.Range("E5").Value = 13.68055134426848
$ws.Range("F5").Value = 31.16299647663767
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 9.852871317558911
$ws.Range("N5").Value = 16.95711894110785
$ws.Range("O5").Value = 23.24408008534252
$ws.Range("B6").Value = 14.03997727294003
$ws.Range("C6").Value = 7.492919034402751
$ws.Range("D6").Value = 9.403955864620448
$ws.Range("E6").Value = 13.67858221507923
$ws.Range("F6").Value = 31.1632521126704
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 9.853328022106322
$ws.Range("N6").Value = 16.95940049030879
$ws.Range("O6").Value = 23.24640902508654
$ws.Range("B7").Value = 14.20107608626508
$ws.Range("C7").Value = 7.645463354524246
$ws.Range("D7").Value = 9.415181565918312
$ws.Range("E7").Value = 13.69250043174787
$ws.Range("F7").Value = 31.16229914720007
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 9.850305757193127
$ws.Range("N7").Value = 16.94373241653498
$ws.Range("O7").Value = 23.23081942634343
$ws.Range("B8").Value = 14.89502703548895
$ws.Range("C8").Value = 8.284357798065301
$ws.Range("D8").Value = 9.467803124235703
$ws.Range("E8").Value = 13.76004961356246
$ws.Range("F8").Value = 31.17842633971458
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 9.840624836818138
$ws.Range("N8").Value = 16.87875304272801
$ws.Range("O8").Value = 23.17631289845744
$ws.Range("B9").Value = 16.18197489636044
$ws.Range("C9").Value = 9.409338431093941
$ws.Range("D9").Value = 9.581590364313151
$ws.Range("E9").Value = 13.91388887860075
$ws.Range("F9").Value = 31.28303896781057
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 9.834764589959894
$ws.Range("N9").Value = 16.76622606063953
$ws.Range("O9").Value = 23.12109789492641
$ws.Range("B10").Value = 17.07266139525742
$ws.Range("C10").Value = 10.15497155153094
$ws.Range("D10").Value = 9.671178845861569
$ws.Range("E10").Value = 14.03928598334532
$ws.Range("F10").Value = 31.40414659051233
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 9.8384386741564
$ws.Range("N10").Value = 16.69258312297936
$ws.Range("O10").Value = 23.11226970770452
$ws.Range("B11").Value = 17.46425112024098
$ws.Range("C11").Value = 10.47621611380621
$ws.Range("D11").Value = 9.713128027480908
$ws.Range("E11").Value = 14.09886259454937
$ws.Range("F11").Value = 31.46875411628483
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 9.841835627596641
$ws.Range("N11").Value = 16.66103151564338
$ws.Range("O11").Value = 23.11517968788406
$ws.Range("B12").Value = 17.61046284851109
$ws.Range("C12").Value = 10.59525922326164
$ws.Range("D12").Value = 9.729175382349311
$ws.Range("E12").Value = 14.12177270271673
$ws.Range("F12").Value = 31.49457606352587
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 9.843369229392481
$ws.Range("N12").Value = 16.64936325708241
$ws.Range("O12").Value = 23.11727895438462
$ws.Range("B13").Value = 17.57906757743619
$ws.Range("C13").Value = 10.56973734848654
$ws.Range("D13").Value = 9.725712244857995
$ws.Range("E13").Value = 14.11682328316196
$ws.Range("F13").Value = 31.48895472476941
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 9.843027960018803
$ws.Range("N13").Value = 16.65186379714669
$ws.Range("O13").Value = 23.11678247114964
$ws.Range("B14").Value = 17.47632226750066
$ws.Range("C14").Value = 10.48606219165941
$ws.Range("D14").Value = 9.714445060651414
$ws.Range("E14").Value = 14.10074049871003
$ws.Range("F14").Value = 31.47085139232831
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 9.841956848300381
$ws.Range("N14").Value = 16.66006595974646
$ws.Range("O14").Value = 23.11533240326936
$ws.Range("B15").Value = 17.41311427131698
$ws.Range("C15").Value = 10.43446885049694
$ws.Range("D15").Value = 9.707564403783337
$ws.Range("E15").Value = 14.09093443517536
$ws.Range("F15").Value = 31.45993884747283
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 9.841332931752449
$ws.Range("N15").Value = 16.6651264224239
$ws.Range("O15").Value = 23.11457410116741
$ws.Range("B16").Value = 17.04678597651366
$ws.Range("C16").Value = 10.1336140753482
$ws.Range("D16").Value = 9.668460588273478
$ws.Range("E16").Value = 14.03544225489192
$ws.Range("F16").Value = 31.40011484482728
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 9.838251325911569
$ws.Range("N16").Value = 16.69468426003274
$ws.Range("O16").Value = 23.11221903929701
$ws.Range("B17").Value = 16.81848428242018
$ws.Range("C17").Value = 9.944433091413247
$ws.Range("D17").Value = 9.644770882132873
$ws.Range("E17").Value = 14.00203821258977
$ws.Range("F17").Value = 31.36584407326631
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 9.836802239698866
$ws.Range("N17").Value = 16.71331576683315
$ws.Range("O17").Value = 23.11254944668219
$ws.Range("B18").Value = 16.6858984137911
$ws.Range("C18").Value = 9.833932860037118
$ws.Range("D18").Value = 9.63125826150848
$ws.Range("E18").Value = 13.98306440106145
$ws.Range("F18").Value = 31.34702887183825
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 9.836131247031185
$ws.Range("N18").Value = 16.72421559409184
$ws.Range("O18").Value = 23.11339136574826
$ws.Range("B19").Value = 16.64079245322498
$ws.Range("C19").Value = 9.796230259898655
$ws.Range("D19").Value = 9.626702832438689
$ws.Range("E19").Value = 13.9766817274431
$ws.Range("F19").Value = 31.34081265722807
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 9.83593199185718
$ws.Range("N19").Value = 16.72793761957631
$ws.Range("O19").Value = 23.11378832951467
$ws.Range("B20").Value = 16.84292003881799
$ws.Range("C20").Value = 9.964746663557429
$ws.Range("D20").Value = 9.647281058889751
$ws.Range("E20").Value = 14.00556946520571
$ws.Range("F20").Value = 31.36939955591586
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 9.836939686961616
$ws.Range("N20").Value = 16.71131342635636
$ws.Range("O20").Value = 23.112446794603
$ws.Range("B21").Value = 17.50655825298763
$ws.Range("C21").Value = 10.51071048128754
$ws.Range("D21").Value = 9.717750184843933
$ws.Range("E21").Value = 14.10545502890242
$ws.Range("F21").Value = 31.47613206948939
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 9.842264757302036
$ws.Range("N21").Value = 16.65764920082
$ws.Range("O21").Value = 23.11573124976854
$ws.Range("B22").Value = 17.92814504276766
$ws.Range("C22").Value = 10.8523370831756
$ws.Range("D22").Value = 9.76474605832472
$ws.Range("E22").Value = 14.17276744050614
$ws.Range("F22").Value = 31.55378763100788
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 9.847185642317953
$ws.Range("N22").Value = 16.62420628187101
$ws.Range("O22").Value = 23.12369099130168
$ws.Range("B23").Value = 17.70428303281447
$ws.Range("C23").Value = 10.67140107951907
$ws.Range("D23").Value = 9.739580718998925
$ws.Range("E23").Value = 14.13666059856419
$ws.Range("F23").Value = 31.51162299396041
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 9.844427785316661
$ws.Range("N23").Value = 16.64190646726256
$ws.Range("O23").Value = 23.11891059629567
$ws.Range("B24").Value = 16.83187677344099
$ws.Range("C24").Value = 9.955568307847191
$ws.Range("D24").Value = 9.646145874248253
$ws.Range("E24").Value = 14.00397226647189
$ws.Range("F24").Value = 31.36778935671545
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 9.836877042042762
$ws.Range("N24").Value = 16.71221809779279
$ws.Range("O24").Value = 23.11249117284317
$ws.Range("B25").Value = 15.84281578287488
$ws.Range("C25").Value = 9.11904071006456
$ws.Range("D25").Value = 9.54972065292208
$ws.Range("E25").Value = 13.87004790491671
$ws.Range("F25").Value = 31.24694186168529
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 9.834946569094271
$ws.Range("N25").Value = 16.79507811521203
$ws.Range("O25").Value = 23.13047551476549
